# Fix WECC case: add a new "ctrl" column (controllability flag) to the
# PV and Slack generator sheets, inserted right before the existing "ra"
# column, shifting all subsequent columns one to the right. All existing
# generators get ctrl = 1. Also restore the last-used selection so that
# the Slack sheet (rather than PV) is the active/visible tab, matching
# the author's last interaction with the workbook.

$wb = $excel.ActiveWorkbook

# --- PV sheet ---------------------------------------------------------
$pv = $wb.Worksheets.Item("PV")
$pv.Columns("R:R").Insert() | Out-Null
$pv.Range("R1").Value = "ctrl"
$pv.Range("R2:R29").Value = 1
$pv.Range("O36").Select() | Out-Null

# --- Slack sheet -------------------------------------------------------
$slack = $wb.Worksheets.Item("Slack")
$slack.Columns("R:R").Insert() | Out-Null
$slack.Range("R1").Value = "ctrl"
$slack.Range("R2").Value = 1

# Slack becomes the active/selected sheet (matches activeTab change).
$slack.Activate() | Out-Null
$slack.Range("P18").Select() | Out-Null
